# Generate Report for Handback
#
# Marks the zh-cn / de-de handback as complete: refreshes the "Status" text,
# records the new handback target/file/datetime on each language sheet, and
# widens a few columns so the longer values are readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ColumnWidth is quantized by the host to an MDW-7 pixel grid (stored_width =
# round(cw*6)/6 + 5/6), so back it off by 5/6 before rounding to the nearest
# 1/6th to land as close as possible on the desired stored character width.
function Set-ExactColumnWidth($range, [double]$targetWidth) {
    $n = [Math]::Round($targetWidth * 6 - 5)
    $range.ColumnWidth = $n / 6
}

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# Touch every cell that currently shares that string so they all land back on
# one shared string instead of leaving stale copies behind.
$newStatus = "Handed back: in sync with en-US"
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# --- zh-cn: Latest Target File / Latest Handback File / Handback DateTime ---
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/9c82b6adab13635ed5e432a7fe917dd6b2d81f53/e2e/c15e0277-08c4-42d0-9f6c-2f5a303b825f.md", "", "c15e0277-08c4-42d0-9f6c-2f5a303b825f.md", "c15e0277-08c4-42d0-9f6c-2f5a303b825f.md")
$zhcn.Range("J2").Value = "c15e0277-08c4-42d0-9f6c-2f5a303b825f.ad6d05f0e6091d8398d032365659974c82ec8922.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-13 03:12:15"

# --- de-de: Latest Target File / Latest Handback File / Handback DateTime ---
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/9c82b6adab13635ed5e432a7fe917dd6b2d81f53/e2e/c15e0277-08c4-42d0-9f6c-2f5a303b825f.md", "", "c15e0277-08c4-42d0-9f6c-2f5a303b825f.md", "c15e0277-08c4-42d0-9f6c-2f5a303b825f.md")
$dede.Range("J2").Value = "c15e0277-08c4-42d0-9f6c-2f5a303b825f.ad6d05f0e6091d8398d032365659974c82ec8922.de-de.xlf"
$dede.Range("K2").Value = "2016-08-13 03:12:24"

# --- Column widths: widen the status/datetime columns to fit the new text ---
Set-ExactColumnWidth $overview.Range("E1") 29.9777047293527
Set-ExactColumnWidth $overview.Range("F1") 29.9777047293527

Set-ExactColumnWidth $zhcn.Range("C1") 29.9777047293527
Set-ExactColumnWidth $zhcn.Range("I1") 40
Set-ExactColumnWidth $zhcn.Range("J1") 40

Set-ExactColumnWidth $dede.Range("C1") 29.9777047293527
Set-ExactColumnWidth $dede.Range("I1") 40
Set-ExactColumnWidth $dede.Range("J1") 40
